# Moved code. No change in functionality
#
# This script updates the two API test rows on the "DataSet" worksheet with
# freshly-generated request/response sample data (new mobile numbers, a new
# Authorization header timestamp/signature), and restores the originally
# selected cells on both worksheets.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("DataSet")
$wsSheet1 = $wb.Worksheets.Item("Sheet1")

# --- Update Body (column L) row 2: refreshed MobileNumber ---
$bodyRow2 = "{""UserId"":""00212029-ba97-468f-b670-b21eb2a93a8e"",`n""EmailAddress"":""info@osn.com"",`n""MobileNumber"":""9651000223"",`n""Packages"":[  `n   3507`n],`n""BirthDate"":null,`n""Address"":null,`n""Country"":null,`n""City"":null,`n""Gender"":null,`n""CustomerUsernameID"":null,`n""Password"":""413703"",`n""Name"":null,`n""Title"":null,`n""LanguagePreference"":null,`n""Email2"":null,`n""Mobile2"":null,`n""Extra"":{  `n   ""MCC"":""965"",`n   ""MNC"":""01"",`n   ""Prod"":""01""`n},`n""CreatedDate"":""2018-12-12T13:00:54.4150898Z"",`n""ExpiryDate"":""2023-12-12T13:00:54.415093Z""}"
$wsData.Range("L2").Value = $bodyRow2

# --- Update Body (column L) row 4: refreshed MobileNumber ---
$bodyRow4 = "{""MobileNumber"" : ""9711000224"", ""EmailAddress"" : ""autoexection@osn.com"", ""Packages"" : [3507], ""Password"" : ""413703"",""extra"": { ""MCC"": ""971"",""MNC"": ""01"",""Prod"": ""01""}}"
$wsData.Range("L4").Value = $bodyRow4

# --- Update Headers (column J) for both request rows (row 2 and row 4) ---
$newHeaders = "Authorization:osnAuth osnauth_x_application_id=6,  osnauth_x_source_id=14, osnauth_x_timestamp=1547128136, osnauth_x_signature=ZjQ5NWEzYTM2ODE5MzJmN2UxOTcwM2UxNmQyMDE4YWU3MGM0MGM4ZjM5YjQxY2VhNTRkMDBhODJlZTdiNjUzMw=="
$wsData.Range("J2").Value = $newHeaders
$wsData.Range("J4").Value = $newHeaders

# --- Restore selections / view state (cells moved/scrolled during editing) ---
$wsSheet1.Select()
$wsSheet1.Range("A1:C6").Select()

$wsData.Select()
$wsData.Range("G1").Select()
$wsData.Range("M2").Select()
